$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value could be mis-parsed as a number/date by Excel's smart
# input parsing (e.g. "1.000" -> 1, "0.6410" -> 0.641, dropping the literal
# formatting that the source data relies on). For those we force the cell to
# Text format first, set the value, then restore the default "Normal" style so
# no stray number-format gets left behind on the cell.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "28.373.69"
$ws.Range("E2").Value = "  +3.85%  "
$ws.Range("D3").Value = "1.802.35"
$ws.Range("E3").Value = "  +1.51%  "
Set-TextValue "D4" "1.002"
$ws.Range("E4").Value = "  +0.11%  "
Set-TextValue "D5" "315.91"
$ws.Range("E5").Value = "  +0.76%  "
Set-TextValue "D6" "1.001"
$ws.Range("E6").Value = "  +0.01%  "
Set-TextValue "D7" "0.5475"
$ws.Range("E7").Value = "  +4.48%  "
Set-TextValue "D8" "0.3852"
$ws.Range("E8").Value = "  +6.40%  "
Set-TextValue "D9" "0.07587"
$ws.Range("E9").Value = "  +2.90%  "
Set-TextValue "D10" "42.59"
$ws.Range("E10").Value = "  +0.14%  "
Set-TextValue "D11" "1.126"
$ws.Range("E11").Value = "  +3.04%  "
Set-TextValue "D12" "1.002"
$ws.Range("E12").Value = "  +0.20%  "
Set-TextValue "D13" "21.21"
$ws.Range("E13").Value = "  +3.31%  "
Set-TextValue "D14" "6.182"
$ws.Range("E14").Value = "  +1.81%  "
Set-TextValue "D15" "7.412"
$ws.Range("E15").Value = "  +6.27%  "
$ws.Range("D16").Value = "1.804.34"
$ws.Range("E16").Value = "  +1.67%  "
Set-TextValue "D17" "92.02"
$ws.Range("E17").Value = "  +4.06%  "
Set-TextValue "D18" "0.00001070"
$ws.Range("E18").Value = "  +2.31%  "
Set-TextValue "D19" "0.06438"
$ws.Range("E19").Value = "  +0.23%  "
Set-TextValue "D20" "1.000"
$ws.Range("E20").Value = "  +0.01%  "
Set-TextValue "D21" "17.31"
$ws.Range("E21").Value = "  +3.24%  "
Set-TextValue "D22" "5.978"
$ws.Range("E22").Value = "  +2.25%  "
$ws.Range("D23").Value = "28.434.39"
$ws.Range("E23").Value = "  +3.79%  "
Set-TextValue "D24" "11.42"
$ws.Range("E24").Value = "  +0.83%  "
Set-TextValue "D25" "2.136"
$ws.Range("E25").Value = "  +3.16%  "
Set-TextValue "D26" "159.33"
$ws.Range("E26").Value = "  +3.52%  "
Set-TextValue "D27" "20.67"
$ws.Range("E27").Value = "  +2.57%  "
Set-TextValue "D28" "2.408"
$ws.Range("E28").Value = "  +2.20%  "
$ws.Range("D29").Value = "2.020.20"
$ws.Range("E29").Value = "  +2.14%  "
Set-TextValue "D30" "123.91"
$ws.Range("E30").Value = "  +2.17%  "
Set-TextValue "D31" "1.118"
$ws.Range("E31").Value = "  +5.14%  "
Set-TextValue "D32" "0.1019"
$ws.Range("E32").Value = "  +4.05%  "
Set-TextValue "D33" "5.767"
$ws.Range("E33").Value = "  +3.57%  "
Set-TextValue "D34" "3.684"
$ws.Range("E34").Value = "  +1.85%  "
Set-TextValue "D35" "0.2319"
$ws.Range("E35").Value = "  +14.48%  "
Set-TextValue "D36" "0.06445"
$ws.Range("E36").Value = "  +7.86%  "
Set-TextValue "D37" "0.02322"
$ws.Range("E37").Value = "  +3.97%  "
Set-TextValue "D38" "8.837"
$ws.Range("E38").Value = "  +9.31%  "
Set-TextValue "D39" "5.151"
$ws.Range("E39").Value = "  +5.92%  "
$ws.Range("E40").Value = "  +3.92%  "
Set-TextValue "D41" "0.6410"
$ws.Range("E41").Value = "  +4.20%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D42" "1.163"
$ws.Range("E42").Value = "  +1.54%  "
$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D43" "1.000"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("E44").Value = "  -3.06%  "
Set-TextValue "D45" "13.50"
$ws.Range("E45").Value = "  +2.96%  "
Set-TextValue "D46" "0.5981"
$ws.Range("E46").Value = "  +3.61%  "
Set-TextValue "D47" "3.683"
$ws.Range("E47").Value = "  +1.45%  "
Set-TextValue "D48" "127.39"
$ws.Range("E48").Value = "  +5.00%  "
Set-TextValue "D49" "1.983"
$ws.Range("E49").Value = "  +4.87%  "
Set-TextValue "D50" "1.149"
$ws.Range("E50").Value = "  +3.30%  "
Set-TextValue "D51" "0.06899"
$ws.Range("E51").Value = "  +2.71%  "
